# Updates cryptos list figures (price + 1h volume change) per the
# Thu May 18 22:40:10 UTC 2023 GitHub Actions refresh, including the two
# coin swaps (rows 41/42 and 49/50 traded places in the ranking).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.910.32"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.810.62"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.52"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4633"
$ws.Range("E7").Value = "  +3.72%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3758"
$ws.Range("E8").Value = "  +0.13%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07485"
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8787"
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.46"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.770.78"
$ws.Range("E12").Value = "  -2.94%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.357"
$ws.Range("E13").Value = "  -0.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.549"
$ws.Range("E14").Value = "  -2.97%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07044"
$ws.Range("E15").Value = "  -0.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.44"
$ws.Range("E16").Value = "  -2.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("E17").Value = "  -0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008758"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.001"
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -2.98%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.916.59"
$ws.Range("E21").Value = "  -1.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.315"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.021.92"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("E25").Value = "  -1.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "151.20"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.60"
$ws.Range("E27").Value = "  -0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.154"
$ws.Range("E28").Value = "  -9.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.310"
$ws.Range("E29").Value = "  -1.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.29"
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7723"
$ws.Range("E32").Value = "  -1.94%  "
$ws.Range("E33").Value = "  -2.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.480"
$ws.Range("E34").Value = "  -0.70%  "
$ws.Range("E35").Value = "  -0.77%  "
$ws.Range("E36").Value = "  -0.07%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01960"
$ws.Range("E38").Value = "  -1.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.443"
$ws.Range("E39").Value = "  +5.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05244"
$ws.Range("E40").Value = "  -1.64%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.262"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5338"
$ws.Range("E42").Value = "  +0.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.908"
$ws.Range("E43").Value = "  +1.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1660"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.593"
$ws.Range("E45").Value = "  -1.71%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5067"
$ws.Range("E46").Value = "  -0.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.31"
$ws.Range("E47").Value = "  -3.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.50"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.668"
$ws.Range("E50").Value = "  -1.89%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06329"
$ws.Range("E51").Value = "  -0.80%  "
